$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.547.29"
$ws.Range("E2").Value = "  -1.19%  "

$ws.Range("D3").Value = "3.502.17"
$ws.Range("E3").Value = "  +0.36%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.91%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.496.73"
$ws.Range("E8").Value = "  +0.23%  "

$ws.Range("E9").Value = "  -2.06%  "

$ws.Range("E10").Value = "  -2.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.90"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").Value = "4.105.00"
$ws.Range("E13").Value = "  +0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.30"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.11%  "

$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").Value = "66.528.77"
$ws.Range("E16").Value = "  -1.22%  "

$ws.Range("E17").Value = "  -2.34%  "

$ws.Range("D18").Value = "3.506.47"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("E19").Value = "  -4.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.70"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.86"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.545"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.83"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -5.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.172"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.10%  "

$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "24.53"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.89"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("E34").Value = "  -5.73%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.25"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.81%  "

$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.29"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "29.69"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +12.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.892"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.60%  "

$ws.Range("E41").Value = "  -5.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.53"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.11%  "

$ws.Range("D44").Value = "2.723.60"
$ws.Range("E44").Value = "  -4.34%  "

$ws.Range("E45").Value = "  -9.86%  "

$ws.Range("E46").Value = "  -2.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.79%  "

$ws.Range("E49").Value = "  -2.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "323.55"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.01"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.10%  "
